# Weekly price-log update: a new week's data is inserted at the top of the
# data block (row 33) and every existing row shifts down by one, with the
# oldest row (61) duplicated into the new last row (62).
#
# Only the columns that actually vary row-to-row need to move:
#   D (Fecha), J (cantidad), K/L/M (precios), P (precio unitario)
# Everything else (A, B, C, E, F, G, H, I, N, O, Q, R) is constant across
# every data row in this sheet, so row 62 can just be copied from row 61.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 33
$lastDataRow  = 61
$newLastRow   = $lastDataRow + 1

# Make row 62 a full copy of row 61 first (constant columns + the moving
# columns, which will hold the correct "shifted" values since row 61 is the
# oldest entry and ends up at the bottom unchanged). Only column D carries a
# non-default (date) style, so only copy NumberFormat there - copying it
# elsewhere would stamp a needless explicit "General" style on every cell.
for ($col = 1; $col -le 18; $col++) {
    $srcCell = $ws.Cells.Item($lastDataRow, $col)
    $dstCell = $ws.Cells.Item($newLastRow, $col)
    if ($col -eq 4) {
        $dstCell.NumberFormat = $srcCell.NumberFormat
    }
    $dstCell.Value = $srcCell.Value()
}

# Shift the moving columns down by one row, from the bottom up so we never
# overwrite a source row before it has been read.
for ($r = $lastDataRow; $r -ge $firstDataRow; $r--) {
    $target = $r + 1

    $srcD = $ws.Cells.Item($r, 4)
    $dstD = $ws.Cells.Item($target, 4)
    $dstD.NumberFormat = $srcD.NumberFormat
    $dstD.Value = $srcD.Value()

    $ws.Cells.Item($target, 10).Value = $ws.Cells.Item($r, 10).Value()   # J
    $ws.Cells.Item($target, 11).Value = $ws.Cells.Item($r, 11).Value()  # K
    $ws.Cells.Item($target, 12).Value = $ws.Cells.Item($r, 12).Value()  # L
    $ws.Cells.Item($target, 13).Value = $ws.Cells.Item($r, 13).Value()  # M
    $ws.Cells.Item($target, 16).Value = $ws.Cells.Item($r, 16).Value()  # P
}

# New, most recent week's entry goes into row 33.
$ws.Cells.Item($firstDataRow, 4).Value = 45062   # Fecha
$ws.Cells.Item($firstDataRow, 10).Value = 25     # cantidad
$ws.Cells.Item($firstDataRow, 11).Value = 12000  # precio minimo
$ws.Cells.Item($firstDataRow, 12).Value = 12000  # precio maximo
$ws.Cells.Item($firstDataRow, 13).Value = 12000  # precio moda/promedio
$ws.Cells.Item($firstDataRow, 16).Value = 1200   # precio unitario
